$wb = $excel.ActiveWorkbook

# Set the last selected cell on the previously-active sheet (Tir_430_50R38)
$wsOld = $wb.Worksheets.Item("Tir_430_50R38")
$wsOld.Activate()
$wsOld.Range("G25").Select() | Out-Null

# Create the new sheet by copying the Tir_213_40R21 template, placed at the end
$src = $wb.Worksheets.Item("Tir_213_40R21")
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$src.Copy($null, $lastSheet) | Out-Null
$ws = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws.Name = "Tir_190_50R10"
$ws.Activate()
$ws.Range("G42").Select() | Out-Null

# Update values for the new tire rig
$ws.Range("H3").Value = "Testrig_Post_190_50R10"
$ws.Range("H5").Value = 0.23241000000000001
$ws.Range("H6").Value = 0.13944000000000001
$ws.Range("H7").Value = 0.17799999999999999
$ws.Range("H9").Value = 9.3000000000000007
$ws.Range("H10").Value = 0
$ws.Range("H11").Value = 0.39100000000000001
$ws.Range("I11").Value = 0.73599999999999999
$ws.Range("H12").Value = 0
$ws.Range("I12").Value = 0
$ws.Range("H13").Value = 87038
$ws.Range("H14").Value = 100
$ws.Range("E15").Value = "Vehicle.Chassis.SuspA1.Linkage.Upright.sWheelCentre"
$ws.Range("E16").Value = "Vehicle.Chassis.Body.sAxleA1"
